$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 463
$ws1.Range("F4").Value = 488
$ws1.Range("F5").Value = 2269
$ws1.Range("F8").Value = 72
$ws1.Range("F9").Value = 1670
$ws1.Range("F10").Value = 1670
$ws1.Range("F13").Value = 1429
$ws1.Range("F16").Value = 625
$ws1.Range("F17").Value = 170
$ws1.Range("F18").Value = 121
$ws1.Range("F19").Value = 7367
$ws1.Range("F20").Value = 8203
$ws1.Range("F29").Value = 254
$ws1.Range("F32").Value = 353
$ws1.Range("F34").Value = 1467
$ws1.Range("F35").Value = 250
$ws1.Range("F38").Value = 296
$ws1.Range("F39").Value = 26
$ws1.Range("F40").Value = 758
$ws1.Range("F44").Value = 253
$ws1.Range("F45").Value = 207
$ws1.Range("F47").Value = 187
$ws1.Range("F48").Value = 174
$ws1.Range("F49").Value = 3

# Sheet: 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 2635
$ws3.Range("F4").Value = 288
$ws3.Range("F6").Value = 17

# Sheet: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 463
$ws4.Range("F8").Value = 488
$ws4.Range("F9").Value = 2269
$ws4.Range("F11").Value = 72
$ws4.Range("F12").Value = 1670
$ws4.Range("F13").Value = 1670
$ws4.Range("F16").Value = 625
$ws4.Range("F18").Value = 170
$ws4.Range("F20").Value = 121
$ws4.Range("F21").Value = 7367
$ws4.Range("F22").Value = 8203
$ws4.Range("F29").Value = 250
$ws4.Range("F33").Value = 296
$ws4.Range("F34").Value = 26
$ws4.Range("F37").Value = 758
$ws4.Range("F44").Value = 253
$ws4.Range("F45").Value = 207
$ws4.Range("F47").Value = 187
